# Swap the contents of rows 3 and 4 (the two species-observation records got
# reordered/re-associated). Only the columns that actually differ between
# the two rows need to be touched; everything else (C, I, K, N, P, S, T, U,
# V, W, Y, Z, AA, AB, AD, AE, AF, AG, AH, AT, AW, AX, AY) is identical
# between the rows and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A","B","D","E","F","G","H","J","M","Q","R","AC","AI","AJ","AK","AO")

foreach ($col in $columns) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")

    $v3 = $cell3.Value2
    $v4 = $cell4.Value2

    $cell3.Value2 = $v4
    $cell4.Value2 = $v3
}
